# Update cryptos list values (Price and Volume(1h) columns) per source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.793.93'
$ws.Range('E2').Value = '  -4.04%  '
$ws.Range('D3').Value = '1.720.07'
$ws.Range('E3').Value = '  -2.56%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'309.21"
$ws.Range('E5').Value = '  -5.84%  '
$ws.Range('D6').Value = "'1.002"
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('D7').Value = "'0.4853"
$ws.Range('E7').Value = '  +3.61%  '
$ws.Range('D8').Value = "'0.3480"
$ws.Range('E8').Value = '  -1.24%  '
$ws.Range('D9').Value = "'42.74"
$ws.Range('E9').Value = '  -2.25%  '
$ws.Range('D10').Value = "'0.07218"
$ws.Range('E10').Value = '  -2.17%  '
$ws.Range('D11').Value = "'1.047"
$ws.Range('E11').Value = '  -3.16%  '
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('D13').Value = "'19.75"
$ws.Range('E13').Value = '  -4.09%  '
$ws.Range('D14').Value = "'5.850"
$ws.Range('E14').Value = '  -2.51%  '
$ws.Range('D15').Value = '1.728.43'
$ws.Range('E15').Value = '  -2.06%  '
$ws.Range('D16').Value = "'6.795"
$ws.Range('E16').Value = '  -5.32%  '
$ws.Range('D17').Value = "'86.31"
$ws.Range('D18').Value = "'0.00001035"
$ws.Range('E18').Value = '  -1.78%  '
$ws.Range('D19').Value = "'0.06404"
$ws.Range('E19').Value = '  -0.22%  '
$ws.Range('E20').Value = '  +0.21%  '
$ws.Range('D21').Value = "'16.50"
$ws.Range('E21').Value = '  -2.44%  '
$ws.Range('D22').Value = "'5.704"
$ws.Range('E22').Value = '  -1.30%  '
$ws.Range('D23').Value = '26.865.81'
$ws.Range('D24').Value = "'10.89"
$ws.Range('E24').Value = '  -2.14%  '
$ws.Range('D25').Value = "'2.052"
$ws.Range('E25').Value = '  -4.85%  '
$ws.Range('D26').Value = "'154.54"
$ws.Range('E26').Value = '  -5.29%  '
$ws.Range('D27').Value = "'19.81"
$ws.Range('E27').Value = '  -0.88%  '
$ws.Range('D28').Value = '1.917.89'
$ws.Range('E28').Value = '  -2.38%  '
$ws.Range('D29').Value = "'2.058"
$ws.Range('E29').Value = '  -5.59%  '
$ws.Range('D30').Value = "'120.22"
$ws.Range('E30').Value = '  -2.12%  '
$ws.Range('D31').Value = "'1.028"
$ws.Range('E31').Value = '  -4.17%  '
$ws.Range('D32').Value = "'0.09270"
$ws.Range('E32').Value = '  -0.32%  '
$ws.Range('E33').Value = '  -2.10%  '
$ws.Range('D34').Value = "'5.345"
$ws.Range('E34').Value = '  -3.50%  '
$ws.Range('D35').Value = "'0.05887"
$ws.Range('E35').Value = '  -3.45%  '
$ws.Range('D36').Value = "'0.02173"
$ws.Range('E36').Value = '  -4.04%  '
$ws.Range('D37').Value = "'1.428"
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('D38').Value = "'10.90"
$ws.Range('E38').Value = '  -6.65%  '
$ws.Range('E39').Value = '  +0.20%  '
$ws.Range('D40').Value = "'0.1979"
$ws.Range('E40').Value = '  -4.25%  '
$ws.Range('E41').Value = '  -3.80%  '
$ws.Range('D42').Value = "'0.5959"
$ws.Range('E42').Value = '  -2.99%  '
$ws.Range('D43').Value = "'1.116"
$ws.Range('D44').Value = "'7.413"
$ws.Range('E44').Value = '  -4.48%  '
$ws.Range('D45').Value = "'12.74"
$ws.Range('E45').Value = '  -2.61%  '
$ws.Range('D46').Value = "'3.574"
$ws.Range('E46').Value = '  -4.43%  '
$ws.Range('D47').Value = "'0.5576"
$ws.Range('E47').Value = '  -3.60%  '
$ws.Range('D48').Value = "'118.82"
$ws.Range('E48').Value = '  -3.84%  '
$ws.Range('D49').Value = "'1.829"
$ws.Range('E49').Value = '  -5.24%  '
$ws.Range('D50').Value = "'0.06623"
$ws.Range('E50').Value = '  -2.77%  '
$ws.Range('D51').Value = "'1.092"
$ws.Range('E51').Value = '  -2.87%  '
